# Add a new row (17) with another DataCamp course entry, matching the
# formatting pattern already used for the preceding rows (13-16): the
# course-name cell (column A) uses the existing "dark navy" course-title
# font/color, and the rating cell (column B) gets its own (new) style slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New course name + rating
$ws.Range("A17").Value = "Introduction to Data Science in Python"
$ws.Range("B17").Value = 4

# A17 should look like A13:A16 (same custom font color used for course titles)
$ws.Range("A17").Font.Color = $ws.Range("A16").Font.Color

# B17 picks up a distinct (new) cell style, as happened in the original edit
# (Excel created a fresh font/style entry for this cell on entry).
$ws.Range("B17").Font.ThemeColor = 1

# Move the active selection down to A18, same as after typing the new row
$ws.Range("A18").Select() | Out-Null
